$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "APS1030592895917216037221674480"
$ws.Range("B2").Value = "DERRIK"
$ws.Range("C2").Value = "HARTMAN"
$ws.Range("D2").Value = "DN014765"
$ws.Range("E2").Value = "GA"

# Row 3
$ws.Range("A3").Value = "APS1030592837760716037221776390"
$ws.Range("B3").Value = "JOCELYN"
$ws.Range("C3").Value = "SHIN"
$ws.Range("D3").Value = "DN014217"
$ws.Range("E3").Value = "GA"

# Row 4 (new)
$ws.Range("A4").Value = "APS1030592736037916037221315742"
$ws.Range("B4").Value = "LYNN"
$ws.Range("C4").Value = "SOLOMON"
$ws.Range("D4").Value = "DNF000361"
$ws.Range("E4").Value = "GA"
